$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1 - italic "help text" row (one cell per column, A1:S1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Set to true if creating a new user to set it as an API service account (this field is ignored for existing users)"
$ws.Range("B1").Value = "Empty values will be ignored. To clear this field, send NONE"
$ws.Range("C1").Value = "Users will be found by Username. Leave a field empty to ignore it."
$ws.Range("D1").Value = "Empty values will be ignored. To clear this field, send NONE"
$ws.Range("E1").Value = "Empty values will be ignored. To clear this field, send NONE"
$ws.Range("F1").Value = "Empty values will be ignored. To clear this field, send NONE"
$ws.Range("G1").Value = "Empty values will be ignored. To clear this field, send NONE"
$ws.Range("H1").Value = "Empty values will be ignored. To clear this field, send NONE"
$ws.Range("I1").Value = "Empty values will be ignored. To clear this field, send NONE"
$ws.Range("J1").Value = "Empty values will be ignored. To clear this field, send NONE"
$ws.Range("K1").Value = "Empty values will be ignored. To clear this field, send NONE"
$ws.Range("L1").Value = "Empty values will be ignored. To clear this field, send NONE"
$ws.Range("M1").Value = "Empty values will be ignored. To clear this field, send NONE"
$ws.Range("N1").Value = "Empty values will be ignored. To clear this field, send NONE"
$ws.Range("O1").Value = "Empty values will be ignored. To clear this field, send NONE"
$ws.Range("P1").Value = "This field is not incremental, the value will be fully replaced. To clear this field, send NONE"
$ws.Range("Q1").Value = "This field is not incremental, the value will be fully replaced. To clear this field, send NONE"
$ws.Range("R1").Value = "This field is not incremental, the value will be fully replaced. To clear this field, send NONE"
$ws.Range("S1").Value = "This field will be set to 'success' if the import is successful, otherwise, it will contain an error message"

# Row 1 uses a new italic font style
$ws.Range("A1:S1").Font.Italic = $true

# ---------------------------------------------------------------------------
# Row 2 - bold header row (A2:S2)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "API Service Account"
$ws.Range("B2").Value = "Active"
$ws.Range("C2").Value = "Username"
$ws.Range("D2").Value = "First Name"
$ws.Range("E2").Value = "Last Name"
$ws.Range("F2").Value = "Email"
$ws.Range("G2").Value = "Phone"
$ws.Range("H2").Value = "Position"
$ws.Range("I2").Value = "Restrict Login Ips"
$ws.Range("J2").Value = "Login Enabled"
$ws.Range("K2").Value = "Custom 1"
$ws.Range("L2").Value = "Custom 2"
$ws.Range("M2").Value = "Custom 3"
$ws.Range("N2").Value = "Custom 4"
$ws.Range("O2").Value = "Custom 5"
$ws.Range("P2").Value = "Teams (not incremental)"
$ws.Range("Q2").Value = "Roles (not incremental)"
$ws.Range("R2").Value = "Teams Managed (not incremental)"
$ws.Range("S2").Value = "Status"

$ws.Range("A2:S2").Font.Bold = $true

# ---------------------------------------------------------------------------
# Placeholder / example rows (4, 6, 8) - blank styled cells for data entry
# ---------------------------------------------------------------------------
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""

$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""

$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""

# Status column example on row 8 keeps its own (Arial grey) style
$ws.Range("Q8").Value = ""

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 94.7109375
$ws.Columns("B").ColumnWidth = 94.7109375
$ws.Columns("C").ColumnWidth = 59.85546875
$ws.Range("D1:E1").EntireColumn.ColumnWidth = 59.85546875
$ws.Range("F1:O1").EntireColumn.ColumnWidth = 53.85546875
$ws.Range("P1:R1").EntireColumn.ColumnWidth = 80.5703125
$ws.Columns("S").ColumnWidth = 90.28515625
